# Make sure subsequent edits are NOT recorded as new tracked changes.
$d = $word.ActiveDocument
$d.TrackRevisions = $false

# Accept the only tracked insertion in the document (the "This "is" the
# first chapter..." paragraph), turning it into plain, non-tracked text.
$d.AcceptAllRevisions()

# Fix the stray curly quotes that were left around "is": the accepted text
# reads  This “is“  the first chapter ...  and should read  This is  the
# first chapter ...
$d.Content.Find.Execute(
    "This " + [char]0x201C + "is" + [char]0x201C + " the first chapter",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This is the first chapter", 2
)

# In case the replace above still produced tracked changes, flatten them too.
$d.AcceptAllRevisions()
